# Adjust logging system configuration
# Apply code style fixes: append a new log row (row 87) to each of the
# four worksheets, mirroring the format of the prior row (86).

$wb = $excel.ActiveWorkbook

# Data for the new row (row 87) per worksheet, in workbook sheet order.
$rowsData = @(
    @{ A = 45873.43682870371; B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x01,0x34"; E = "0x14"; F = 380; G = 759863127514710945038336.0; H = 308; I = 14 },
    @{ A = 45873.43682870371; B = "0x01,0x7c"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x01,0x34"; E = "0xe";  F = 380; G = 568432987514711010443264.0;  H = 308; I = 14 },
    @{ A = 45873.43682870371; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x77"; E = "0x7";  F = 130; G = 568631262647113970876416.0;  H = 119; I = 7 },
    @{ A = 45873.43682870371; B = "0x00,0x82"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x75"; E = "0x3";  F = 130; G = 985046333984776009023488.0;  H = 117; I = 3 }
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $data = $rowsData[$i]

    $newRow = 87
    $prevRow = 86

    # Match the number format of the prior row's timestamp cell so the new
    # cell renders the same way (date/time format). The other columns stay
    # on the default "General" format, same as the rest of the table.
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
